# chl-transect-info.xlsx — minor edits to metadata templates XLS
#
# 1. ColumnHeaders!G5 (missingValueCodeExplanation for bottle_other_method):
#      "No associated bottle from the other method"
#      -> "No associated niskin from the other method"
#    ColumnHeaders!F5,F14,F15,F17,F18,F19 (missingValueCode): "NA" -> "NaN"
#
# 2. Personnel!G6, G7, G8 (role for the three creators who are PIs):
#      "Principal Investigator" -> "PI"
#
# 3. View-state changes: the active/selected sheet moves from "Personnel"
#    back to "ColumnHeaders", with new active-cell selections on each.

$wb = $excel.ActiveWorkbook

$colHeaders = $wb.Worksheets.Item("ColumnHeaders")
$personnel = $wb.Worksheets.Item("Personnel")

# --- Personnel sheet: "Principal Investigator" role -> "PI" --------------
# (written first so the new shared-string order matches the source edit:
#  PI, then NaN, then the niskin explanation)
$personnel.Range("G6").Value = "PI"
$personnel.Range("G7").Value = "PI"
$personnel.Range("G8").Value = "PI"

# --- ColumnHeaders sheet ---------------------------------------------------
# missingValueCode "NA" -> "NaN" for every row that used it
$colHeaders.Range("F5").Value = "NaN"
$colHeaders.Range("F14").Value = "NaN"
$colHeaders.Range("F15").Value = "NaN"
$colHeaders.Range("F17").Value = "NaN"
$colHeaders.Range("F18").Value = "NaN"
$colHeaders.Range("F19").Value = "NaN"

# missingValueCodeExplanation for bottle_other_method: "bottle" -> "niskin"
$colHeaders.Range("G5").Value = "No associated niskin from the other method"

# --- View state: active tab moves from Personnel to ColumnHeaders --------
# Update Personnel's remembered selection first (this temporarily activates
# it), then activate ColumnHeaders last so it ends up the visible/active
# sheet with its own selection, matching the saved workbook view state.
$personnel.Range("C5").Select()

$colHeaders.Activate()
$colHeaders.Range("G6").Select()
